# Update result values for result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.2802
$ws.Range("A9").Value = -21.9371
$ws.Range("C12").Value = -11.1312
$ws.Range("E15").Value = 16.2242
$ws.Range("A18").Value = -22.21990000000001
$ws.Range("A20").Value = -21.06459999999998
$ws.Range("C26").Value = -12.68720000000001
$ws.Range("A27").Value = -21.78019999999999
$ws.Range("C27").Value = -12.59979999999999
$ws.Range("C29").Value = -11.51110000000001
$ws.Range("C37").Value = -13.8248
$ws.Range("C38").Value = -13.1208
$ws.Range("E38").Value = 16.21249999999999
$ws.Range("E44").Value = 16.79319999999998
$ws.Range("C51").Value = -12.01499999999999
$ws.Range("E51").Value = 17.34010000000001
$ws.Range("C55").Value = -13.9064
$ws.Range("E57").Value = 16.6638
$ws.Range("E63").Value = 18.48610000000002
$ws.Range("A69").Value = -21.9132
$ws.Range("C69").Value = -11.9161
$ws.Range("C70").Value = -11.9293
$ws.Range("E70").Value = 17.49070000000001
$ws.Range("A76").Value = -20.33839999999998
$ws.Range("A82").Value = -21.96579999999999
$ws.Range("C83").Value = -13.55169999999999
$ws.Range("E99").Value = 16.4456
$ws.Range("C102").Value = -13.4282

$wb.Save()
